# Reconfigured base project to use interrupts
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = "HV_enable_OUT"

$ws.Range("E7").Value = "EXTI2_Select0"
$ws.Range("F7").Value = "GM_pulse_IRQ"

$ws.Range("E10").Value = "EXTI4_Select2"
$ws.Range("F10").Value = "POWER_button_IRQ"

$ws.Range("F11").Value = "POWER_enable_OUT"

$ws.Range("E12").Value = "EXTI6_Select0"
$ws.Range("F12").Value = "USB_sense_IRQ"

$ws.Range("E24").Value = "EXTI10_Select2"
$ws.Range("F24").Value = "MENU_button_IRQ"

$ws.Range("F25").Value = "LED1_red_OUT"

$ws.Range("F26").Value = "LED2_green_OUT"
